$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 3
$ws.Range("H3").Value = 51931.4
$ws.Range("J3").Value = 51931.4
$ws.Range("L3").Value = 51931.4
$ws.Range("N3").Value = -52159.4
# Row 17
$ws.Range("H17").Value = 197020.28
$ws.Range("J17").Value = 197020.28
$ws.Range("L17").Value = 591060.84
$ws.Range("N17").Value = -591396.84
# Row 52
$ws.Range("H52").Value = 46323.637
$ws.Range("J52").Value = 48505.715
$ws.Range("L52").Value = 145517.145
$ws.Range("N52").Value = -145837.145
# Row 102
$ws.Range("H102").Value = 51931.4
$ws.Range("J102").Value = 51931.4
$ws.Range("L102").Value = 51931.4
$ws.Range("N102").Value = -58421.4
# Row 125
$ws.Range("H125").Value = 50618
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 50618
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 455562
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -460482
# Row 132
$ws.Range("H132").Value = 1978.7778
$ws.Range("I132").Value = 1742.4722
$ws.Range("J132").Value = 2924
$ws.Range("K132").Value = 5227.4166
$ws.Range("L132").Value = 8772
$ws.Range("M132").Value = -2697.4166
$ws.Range("N132").Value = -13832
# Row 138
$ws.Range("H138").Value = 3625.114
$ws.Range("I138").Value = 1289.3793
$ws.Range("J138").Value = 4979.84
$ws.Range("K138").Value = 3868.1379
$ws.Range("L138").Value = 14939.52
$ws.Range("M138").Value = 1271.8621
$ws.Range("N138").Value = -25219.52

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 23
$ws.Range("H23").Value = 41502.6
$ws.Range("J23").Value = 31876.75
$ws.Range("L23").Value = 31876.75
$ws.Range("N23").Value = -32394.75
# Row 44
$ws.Range("H44").Value = 333368320
$ws.Range("J44").Value = 333368320
$ws.Range("L44").Value = 333368320
$ws.Range("N44").Value = -333369296
# Row 55
$ws.Range("H55").Value = 59000
$ws.Range("J55").Value = 59000
$ws.Range("L55").Value = 59000
$ws.Range("N55").Value = -59630
# Row 62
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
# Row 65
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
# Row 80
$ws.Range("H80").Value = 250016880
$ws.Range("I80").Value = 22500
$ws.Range("J80").Value = 1000000000
$ws.Range("K80").Value = 22500
$ws.Range("L80").Value = 1000000000
$ws.Range("M80").Value = -21502
$ws.Range("N80").Value = -1000001996
# Row 83
$ws.Range("H83").Value = 250016880
$ws.Range("I83").Value = 22500
$ws.Range("J83").Value = 1000000000
$ws.Range("K83").Value = 67500
$ws.Range("L83").Value = 3000000000
$ws.Range("M83").Value = -62508
$ws.Range("N83").Value = -3000009984
# Row 122
$ws.Range("H122").Value = 5683610.5
$ws.Range("I122").Value = 1755.0588
$ws.Range("K122").Value = 5265.1764
$ws.Range("M122").Value = -2815.1764
# Row 132
$ws.Range("H132").Value = 4222.6963
$ws.Range("I132").Value = 1589.0646
$ws.Range("J132").Value = 7488.4
$ws.Range("K132").Value = 4767.1938
$ws.Range("L132").Value = 22465.2
$ws.Range("M132").Value = -2237.1938
$ws.Range("N132").Value = -27525.2

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 36
$ws.Range("H36").Value = 9999.799999999999
$ws.Range("I36").Value = 3499.75
$ws.Range("J36").Value = 36000
$ws.Range("K36").Value = 3499.75
$ws.Range("L36").Value = 36000
$ws.Range("M36").Value = -2965.75
$ws.Range("N36").Value = -37068
# Row 100
$ws.Range("H100").Value = 33330
$ws.Range("J100").Value = 33330
$ws.Range("L100").Value = 33330
$ws.Range("N100").Value = -35494

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 43
$ws.Range("H43").Value = 33750
$ws.Range("J43").Value = 33750
$ws.Range("L43").Value = 33750
$ws.Range("N43").Value = -34118
# Row 50
$ws.Range("H50").Value = 22129.166
$ws.Range("J50").Value = 22129.166
$ws.Range("L50").Value = 22129.166
$ws.Range("N50").Value = -23379.166
# Row 55
$ws.Range("H55").Value = 14443.777
$ws.Range("I55").Value = 10000
$ws.Range("J55").Value = 14999.25
$ws.Range("K55").Value = 10000
$ws.Range("L55").Value = 14999.25
$ws.Range("M55").Value = -9685
$ws.Range("N55").Value = -15629.25
# Row 57
$ws.Range("H57").Value = 23040.666
$ws.Range("J57").Value = 24420.75
$ws.Range("L57").Value = 24420.75
$ws.Range("N57").Value = -25540.75
# Row 101
$ws.Range("H101").Value = 33750
$ws.Range("J101").Value = 33750
$ws.Range("L101").Value = 33750
$ws.Range("N101").Value = -40240

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 36
$ws.Range("H36").Value = 1902.5
$ws.Range("I36").Value = 1850
$ws.Range("J36").Value = 1920
$ws.Range("K36").Value = 5550
$ws.Range("L36").Value = 5760
$ws.Range("M36").Value = -5381
$ws.Range("N36").Value = -6098
# Row 39
$ws.Range("H39").Value = 10147.619
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 10147.619
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 30442.857
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -31030.857
# Row 55
$ws.Range("H55").Value = 3000
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
# Row 137
$ws.Range("H137").Value = 55558484
$ws.Range("I137").Value = 62502544
$ws.Range("J137").Value = 6000
$ws.Range("K137").Value = 187507632
$ws.Range("L137").Value = 18000
$ws.Range("M137").Value = -187502532
$ws.Range("N137").Value = -28200

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 57
$ws.Range("H57").Value = 20000
$ws.Range("J57").Value = 20000
$ws.Range("L57").Value = 20000
$ws.Range("N57").Value = -21640
# Row 80
$ws.Range("H80").Value = 5374.95
$ws.Range("I80").Value = 8399.875
$ws.Range("J80").Value = 3358.3333
$ws.Range("K80").Value = 8399.875
$ws.Range("L80").Value = 3358.3333
$ws.Range("M80").Value = -7401.875
$ws.Range("N80").Value = -5354.3333
# Row 83
$ws.Range("H83").Value = 5374.95
$ws.Range("I83").Value = 8399.875
$ws.Range("J83").Value = 3358.3333
$ws.Range("K83").Value = 41999.375
$ws.Range("L83").Value = 16791.6665
$ws.Range("M83").Value = -37007.375
$ws.Range("N83").Value = -26775.6665
# Row 126
$ws.Range("H126").Value = 2922.3635
$ws.Range("I126").Value = 1964.9166
$ws.Range("J126").Value = 4071.3
$ws.Range("K126").Value = 5894.7498
$ws.Range("L126").Value = 12213.9
$ws.Range("M126").Value = -3424.7498
$ws.Range("N126").Value = -17153.9

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 5277
$ws.Range("I7").Value = 3620.3
$ws.Range("J7").Value = 7117.778
$ws.Range("K7").Value = 3620.3
$ws.Range("L7").Value = 7117.778
$ws.Range("M7").Value = -3508.3
$ws.Range("N7").Value = -7341.778
# Row 20
$ws.Range("H20").Value = 12000
$ws.Range("J20").Value = 12000
$ws.Range("L20").Value = 12000
$ws.Range("N20").Value = -12452
# Row 40
$ws.Range("H40").Value = 4222.0303
$ws.Range("I40").Value = 3762.7917
$ws.Range("K40").Value = 3762.7917
$ws.Range("M40").Value = -3626.7917
# Row 46
$ws.Range("H46").Value = 780.82355
$ws.Range("I46").Value = 653.4286
$ws.Range("J46").Value = 870
$ws.Range("K46").Value = 653.4286
$ws.Range("L46").Value = 870
$ws.Range("M46").Value = -465.4286
$ws.Range("N46").Value = -1246
# Row 122
$ws.Range("H122").Value = 4974.8125
$ws.Range("I122").Value = 3923.946
$ws.Range("J122").Value = 8509.546
$ws.Range("K122").Value = 11771.838
$ws.Range("L122").Value = 25528.638
$ws.Range("M122").Value = -9321.838
$ws.Range("N122").Value = -30428.638
# Row 126
$ws.Range("H126").Value = 5277
$ws.Range("I126").Value = 3620.3
$ws.Range("J126").Value = 7117.778
$ws.Range("K126").Value = 10860.9
$ws.Range("L126").Value = 21353.334
$ws.Range("M126").Value = -8390.900000000001
$ws.Range("N126").Value = -26293.334

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 92
$ws.Range("H92").Value = 40550
$ws.Range("J92").Value = 40550
$ws.Range("L92").Value = 40550
$ws.Range("N92").Value = -45542
